$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39: record hours worked (6) in column B. Column D text is unchanged.
$ws.Range("B39").Value = 6

# Row 40: new time-log entry (previously a blank placeholder row).
$ws.Range("A40").Value = 43535
$ws.Range("B40").Value = 1.5
$ws.Range("D40").Value = "Indie Project/Week 6: worked out issues concerning access to AWS EC2 instance.  Project is now deployed and public, and SSH works from classroom."
$ws.Rows.Item(40).RowHeight = 30

# Row 42: remove the stray "6:45 - x" placeholder note.
$ws.Range("D42").Clear()

# Reflect the author's final cursor position/scroll in the saved view.
$null = $ws.Range("D41").Select()
